$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2108626198083067
$ws.Range("C2").Value = 0.5175718849840255
$ws.Range("J2").Value = 0.01277955271565495
$ws.Range("P2").Value = 0.1405750798722045
$ws.Range("S2").Value = 0.1182108626198083
$ws.Range("B3").Value = 0.005882352941176471
$ws.Range("C3").Value = 0.02941176470588235
$ws.Range("J3").Value = 0.02941176470588235
$ws.Range("P3").Value = 0.7705882352941177
$ws.Range("S3").Value = 0.1647058823529412
$ws.Range("J4").Value = 0.04838709677419355
$ws.Range("P4").Value = 0.7580645161290323
$ws.Range("S4").Value = 0.1935483870967742
$ws.Range("B6").Value = 0.06504065040650407
$ws.Range("D6").Value = 0.01219512195121951
$ws.Range("F6").Value = 0.04878048780487805
$ws.Range("J6").Value = 0.2154471544715447
$ws.Range("O6").Value = 0.01626016260162602
$ws.Range("Q6").Value = 0.1869918699186992
$ws.Range("R6").Value = 0.08130081300813008
$ws.Range("S6").Value = 0.3739837398373984
$ws.Range("B7").Value = 0.136150234741784
$ws.Range("D7").Value = 0.02347417840375587
$ws.Range("F7").Value = 0.04225352112676056
$ws.Range("J7").Value = 0.1408450704225352
$ws.Range("O7").Value = 0.03755868544600939
$ws.Range("Q7").Value = 0.1267605633802817
$ws.Range("R7").Value = 0.0892018779342723
$ws.Range("S7").Value = 0.4037558685446009
$ws.Range("B8").Value = 0.1077981651376147
$ws.Range("D8").Value = 0.02981651376146789
$ws.Range("F8").Value = 0.07339449541284404
$ws.Range("J8").Value = 0.1192660550458716
$ws.Range("O8").Value = 0.01834862385321101
$ws.Range("Q8").Value = 0.1513761467889908
$ws.Range("R8").Value = 0.1123853211009174
$ws.Range("S8").Value = 0.3876146788990826
$ws.Range("B9").Value = 0.08383233532934131
$ws.Range("D9").Value = 0.01796407185628742
$ws.Range("E9").Value = 0.005988023952095809
$ws.Range("F9").Value = 0.09580838323353294
$ws.Range("J9").Value = 0.1137724550898204
$ws.Range("O9").Value = 0.02395209580838323
$ws.Range("Q9").Value = 0.1796407185628743
$ws.Range("R9").Value = 0.0658682634730539
$ws.Range("S9").Value = 0.4131736526946108
$ws.Range("B10").Value = 0.1041347626339969
$ws.Range("D10").Value = 0.02909647779479326
$ws.Range("F10").Value = 0.09264931087289434
$ws.Range("J10").Value = 0.108728943338438
$ws.Range("O10").Value = 0.03139356814701378
$ws.Range("Q10").Value = 0.1753445635528331
$ws.Range("R10").Value = 0.07963246554364471
$ws.Range("S10").Value = 0.3790199081163859
$ws.Range("G11").Value = 0.1424418604651163
$ws.Range("J11").Value = 0.09593023255813954
$ws.Range("K11").Value = 0.2034883720930233
$ws.Range("L11").Value = 0.5494186046511628
$ws.Range("S11").Value = 0.008720930232558139
$ws.Range("G12").Value = 0.708994708994709
$ws.Range("J12").Value = 0.2222222222222222
$ws.Range("K12").Value = 0.02116402116402116
$ws.Range("L12").Value = 0.01587301587301587
$ws.Range("S12").Value = 0.03174603174603174
$ws.Range("G13").Value = 0.68
$ws.Range("J13").Value = 0.22
$ws.Range("S13").Value = 0.1
$ws.Range("F15").Value = 0.01219512195121951
$ws.Range("H15").Value = 0.1016260162601626
$ws.Range("I15").Value = 0.08943089430894309
$ws.Range("J15").Value = 0.2926829268292683
$ws.Range("K15").Value = 0.05284552845528456
$ws.Range("M15").Value = 0.01626016260162602
$ws.Range("N15").Value = 0.004065040650406504
$ws.Range("O15").Value = 0.05691056910569105
$ws.Range("S15").Value = 0.3739837398373984
$ws.Range("F16").Value = 0.009523809523809525
$ws.Range("H16").Value = 0.1476190476190476
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.3761904761904762
$ws.Range("K16").Value = 0.1095238095238095
$ws.Range("M16").Value = 0.02380952380952381
$ws.Range("O16").Value = 0.05238095238095238
$ws.Range("S16").Value = 0.2095238095238095
$ws.Range("F17").Value = 0.01265822784810127
$ws.Range("H17").Value = 0.1594936708860759
$ws.Range("I17").Value = 0.0810126582278481
$ws.Range("J17").Value = 0.4405063291139241
$ws.Range("K17").Value = 0.0759493670886076
$ws.Range("M17").Value = 0.02025316455696203
$ws.Range("N17").Value = 0.002531645569620253
$ws.Range("O17").Value = 0.04556962025316456
$ws.Range("S17").Value = 0.1620253164556962
$ws.Range("F18").Value = 0.01477832512315271
$ws.Range("H18").Value = 0.2167487684729064
$ws.Range("I18").Value = 0.0541871921182266
$ws.Range("J18").Value = 0.3891625615763547
$ws.Range("K18").Value = 0.08374384236453201
$ws.Range("M18").Value = 0.009852216748768473
$ws.Range("O18").Value = 0.06403940886699508
$ws.Range("S18").Value = 0.167487684729064
$ws.Range("F19").Value = 0.01298701298701299
$ws.Range("H19").Value = 0.1998556998556998
$ws.Range("I19").Value = 0.06277056277056277
$ws.Range("J19").Value = 0.3751803751803752
$ws.Range("K19").Value = 0.1334776334776335
$ws.Range("M19").Value = 0.02453102453102453
$ws.Range("N19").Value = 0.001443001443001443
$ws.Range("O19").Value = 0.06926406926406926
$ws.Range("S19").Value = 0.1204906204906205
